$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.083.58'
$ws.Range('D3').Value = '1.639.93'
$ws.Range('E3').Value = '  -1.63%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.41'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5215'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.20%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.002'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2602'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.51%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06309'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.19%  '
$ws.Range('E10').Value = '  -1.92%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07688'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.22%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.412'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.35%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.623.24'
$ws.Range('E13').Value = '  -3.08%  '
$ws.Range('D14').Value = '1.863.25'
$ws.Range('E14').Value = '  -1.74%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5557'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.13%  '
$ws.Range('D16').Value = '0.0₅8189'
$ws.Range('E16').Value = '  +3.18%  '
$ws.Range('E17').Value = '  -1.92%  '
$ws.Range('D18').Value = '26.083.56'
$ws.Range('E18').Value = '  -0.27%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.002'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.719'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '189.25'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.62%  '
$ws.Range('E22').Value = '  -0.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.170'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.002'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.54'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.49%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.417'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.94%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1203'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.38%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.85'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.394'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05890'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -7.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.255'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.441'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.55%  '
$ws.Range('E33').Value = '  -0.18%  '
$ws.Range('E34').Value = '  +0.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9833'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.97%  '
$ws.Range('B36').Value = 'MXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.763'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.60%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.393'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.64%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5651'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01616'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8536'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.60%  '
$ws.Range('E41').Value = '  -0.25%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.704'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.60%  '
$ws.Range('D43').Value = '1.026.79'
$ws.Range('E43').Value = '  -7.43%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.09'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').Value = '1.790.22'
$ws.Range('E45').Value = '  -1.60%  '
$ws.Range('E46').Value = '  -1.28%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.78'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.84%  '
$ws.Range('E48').Value = '  +0.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.059'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.38%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05151'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4218'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.64%  '
